$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.916.81"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.634.69"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "211.88"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.92%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.523"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.18%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "23.16"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  -3.39%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0612"
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0879"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.866.35"
$ws.Range("D13").Value = "1.635.09"
$ws.Range("E13").Value = "  -1.01%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.06"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -0.69%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.566"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.35%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "65.11"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "27.921.08"
$ws.Range("E17").Value = "  -0.23%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "229.99"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  -3.47%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -3.75%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "152.96"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +0.34%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "6.94"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.50%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "15.62"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.06"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.396.28"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  +9.77%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  +0.00%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.559"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -1.10%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.869"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -0.13%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "66.83"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -3.69%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.51"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "1.775.82"
$ws.Range("E47").Value = "  -1.03%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "87.59"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  -0.28%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "7.52"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -2.67%  "
